$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the formatting of the
# existing header cells (e.g. G1 - bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Save column values for rows 2-7
$saveValues = @(1, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
